$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C. Excel automatically shifts the
# "Sex"/"Time" columns right and extends the A1:D1 / A2:D2 merges to
# A1:E1 / A2:E2, copying column B's formatting into the new column C.
$ws.Range("C1").EntireColumn.Insert()

# The old single "Name" column (B) becomes "First Name"; the new column
# (C) becomes "Last Name", splitting the name for sortability.
$ws.Range("B3").Value = "First Name"
$ws.Range("C3").Value = "Last Name"

# Resize the name columns to fit the new headers, and nudge the "Time"
# column back to its intended width.
$ws.Columns.Item(2).ColumnWidth = 20.0
$ws.Columns.Item(3).ColumnWidth = 21.33
$ws.Columns.Item(5).ColumnWidth = 12.5

# Restore the selection as recorded after the edit.
$ws.Range("E11").Select() | Out-Null
